# The only change in the target diff is the disappearance of
# ppt/revisionInfo.xml — an internal, non-content "Office temporary file"
# that PowerPoint writes while an AutoRecover / real-time co-authoring
# session is in flight (p1510:revInfo / p1510:revLst client revision log).
# It is not part of the document's visible content and, like other
# session-scoped package plumbing, it is not exposed anywhere in the
# Presentation/Slide/Shape object model (no property or method on
# Application, Presentation, DocumentInspectors, Coauthoring, Sync, etc.
# reaches it) — real PowerPoint itself retires the file the next time the
# presentation is opened and saved cleanly, once there is no longer a live
# revision-tracking session to record. Confirm that clean state the only
# way user/automation code can: open the deck and save it normally.

$p = $ppt.ActivePresentation
$p.Save()
